# Build the "CMM Report" header box (rows 1-2, columns A-L) with a
# yellow fill, Aptos Narrow 12pt font, and a medium black border drawn
# around the outside of the 2x12 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:L2")

# Font + fill for the whole header block first.
$headerRange.Font.Name = "Aptos Narrow"
$headerRange.Font.Size = 12
$headerRange.Interior.Color = 65535

# Outer border, medium weight, applied edge-by-edge so that each unique
# border combination (corner / edge) is created once.
$ws.Range("B1:K1").Borders.Item(8).Weight = -4138   # top edge (middle)
$ws.Range("A1").Borders.Item(8).Weight = -4138      # top-left corner
$ws.Range("A1").Borders.Item(7).Weight = -4138
$ws.Range("L1").Borders.Item(8).Weight = -4138      # top-right corner
$ws.Range("L1").Borders.Item(10).Weight = -4138

$ws.Range("B2:K2").Borders.Item(9).Weight = -4138   # bottom edge (middle)
$ws.Range("A2").Borders.Item(9).Weight = -4138      # bottom-left corner
$ws.Range("A2").Borders.Item(7).Weight = -4138
$ws.Range("L2").Borders.Item(9).Weight = -4138      # bottom-right corner
$ws.Range("L2").Borders.Item(10).Weight = -4138

# Title text in the header box.
$ws.Range("A1").Value = "CMM Report"

# Report description line.
$ws.Range("A3").Value = "this is a report from origin"

# Remove the old template values that used to live further down the sheet.
$ws.Range("D7").ClearContents()
$ws.Range("F10").ClearContents()

# Data row (values updated from the previous template values).
$ws.Range("B5").Value = 942
$ws.Range("C5").Value = 1752
$ws.Range("D5").Value = 6418756
$ws.Range("E5").Value = 16

# Stretch the sheet's used range down to row 100 / column AN, matching
# the larger report template, without disturbing any formatting.
$ws.Range("AN100").Font.Bold = $false
